# Legs Update and Sesi 1 Update!
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- "Legs" data updates (row 6 = L3, row 8 = R2, row 9 = R3) ---

# L3 leg: alpha_90 input (G6) updated; dependent formulas (H6, K6) recalc automatically.
$ws.Range("G6").Value = 2120

# R2 leg: beta_0 input (D8) and alpha_90 input (G8) updated; H8 recalcs automatically.
$ws.Range("D8").Value = 1270
$ws.Range("G8").Value = 750

# R3 leg: beta_0 (C9) and alpha_0 (D9) updated; H9, I9, K9, L9 recalc automatically.
$ws.Range("C9").Value = 1400
$ws.Range("D9").Value = 1250

# --- Selection moved from F14 to C14 (Sesi 1 update) ---
$ws.Range("C14").Select()
